$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings: issue number and reporting week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/6/2024  Through  5/12/2024"

# --- Helper cells (far outside the used range) holding genuine TEXT values "0" and "***.*" ---
# Using a formula that evaluates to a text string guarantees the result is stored as a
# shared-string (t="s") rather than being auto-coerced to a number by plain Value assignment.
$ws.Range("ZZ1").Formula = "=""0"""
$ws.Range("ZZ2").Formula = "=""***.*"""

# --- Numeric value updates across the crime-statistics table (rows 15-33) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 50
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 9
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 64
$ws.Range("K16").Value = 31.25
$ws.Range("L16").Value = 40
$ws.Range("M16").Value = 1.204819277108
$ws.Range("N16").Value = -48.148148148148
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 8.695652173913
$ws.Range("I17").Value = 113
$ws.Range("J17").Value = 105
$ws.Range("K17").Value = 7.619047619047
$ws.Range("L17").Value = 34.523809523809
$ws.Range("M17").Value = 76.5625
$ws.Range("N17").Value = 43.037974683544
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 56
$ws.Range("K18").Value = -7.142857142857
$ws.Range("L18").Value = 30
$ws.Range("M18").Value = -43.478260869565
$ws.Range("N18").Value = -83.544303797468
$ws.Range("C19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -33.898305084745
$ws.Range("I19").Value = 219
$ws.Range("J19").Value = 241
$ws.Range("K19").Value = -9.128630705394
$ws.Range("L19").Value = 12.307692307692
$ws.Range("M19").Value = 51.03448275862
$ws.Range("N19").Value = 44.078947368421
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -15.384615384615
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 53
$ws.Range("H20").Value = -43.396226415094
$ws.Range("I20").Value = 143
$ws.Range("J20").Value = 192
$ws.Range("K20").Value = -25.520833333333
$ws.Range("L20").Value = 57.142857142857
$ws.Range("M20").Value = 134.426229508197
$ws.Range("N20").Value = -79.658605974395
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -25.641025641025
$ws.Range("F21").Value = 114
$ws.Range("G21").Value = 160
$ws.Range("H21").Value = -28.75
$ws.Range("I21").Value = 622
$ws.Range("J21").Value = 664
$ws.Range("K21").Value = -6.325301204819
$ws.Range("L21").Value = 30.672268907563
$ws.Range("M21").Value = 37.610619469026
$ws.Range("N21").Value = -56.289529163738
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("I23").Value = 23
$ws.Range("J23").Value = 21
$ws.Range("K23").Value = 9.523809523809
$ws.Range("L23").Value = 76.923076923076
$ws.Range("M23").Value = 76.923076923076
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 119
$ws.Range("G24").Value = 98
$ws.Range("H24").Value = 21.428571428571
$ws.Range("I24").Value = 515
$ws.Range("J24").Value = 492
$ws.Range("K24").Value = 4.674796747967
$ws.Range("L24").Value = 27.475247524752
$ws.Range("M24").Value = 10.042735042735
$ws.Range("C25").Value = 24
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 78
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 13.043478260869
$ws.Range("I25").Value = 302
$ws.Range("J25").Value = 279
$ws.Range("K25").Value = 8.243727598566
$ws.Range("L25").Value = 48.768472906403
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -46.153846153846
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 6.976744186046
$ws.Range("I26").Value = 200
$ws.Range("J26").Value = 186
$ws.Range("K26").Value = 7.52688172043
$ws.Range("L26").Value = 21.951219512195
$ws.Range("M26").Value = 41.843971631205
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("L27").Value = 0
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 60
$ws.Range("I28").Value = 30
$ws.Range("K28").Value = 100
$ws.Range("L28").Value = 50
$ws.Range("L33").Value = -50

# --- Cells that change from a numeric value to the textual placeholder "0" or "***.*" ---
# First copy the (General-format) style already used by neighboring text cells, then paste
# only the value from the helper cell so the destination keeps style 14 but becomes text.
$ws.Range("D14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C14").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("ZZ2").Copy()
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C27").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("ZZ2").Copy()
$ws.Range("E27").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("ZZ2").Copy()
$ws.Range("E28").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C29").PasteSpecial(-4163)
$ws.Range("D14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("ZZ1").Copy()
$ws.Range("C30").PasteSpecial(-4163)

# --- Remove the temporary helper cells ---
$ws.Range("ZZ1").Clear()
$ws.Range("ZZ2").Clear()
